$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append below the existing data (through 26/04 -> row 238)
$data = @(
    @(44308, 1, 6, 244.798041615667),
    @(44309, 0, 6, 244.798041615667),
    @(44310, 0, 6, 244.798041615667),
    @(44311, 0, 5, 203.9983680130559),
    @(44312, 0, 5, 203.9983680130559)
)

$startRow = 234
$lastExistingRow = 233

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy formatting (style/number format) from the last existing data row
    # so the new rows look identical to the ones above them.
    $ws.Range("A$lastExistingRow").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = 0
